$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update wheel-radius parameter (C15: 15.6 -> 15.65)
$ws.Range("C15").Value = 15.65

# Update the C5 formula: use 180 (instead of 120) and PI() (instead of 3.14)
$ws.Range("C5").Formula = "=180*D5*2/(2*PI()*15.65) /100"

# D10 gets the "ms" label
$ws.Range("D10").Value = "ms"

# New notes / scratch calculations added lower on the sheet
$ws.Range("C47").Value = 42717
$ws.Range("C47").NumberFormat = "d-mmm"

$ws.Range("C56").Formula = "=13/12"
$ws.Range("C57").Formula = "=530/2/180 * 2*PI()*C15 /120"
$ws.Range("C59").Formula = "=130/(2*PI()*C15) * 180"
$ws.Range("C60").Formula = "=28000/(2*PI()*15.65)"

# Update sheet view: scrolled/selection position
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("C57").Select()
